$wb = $excel.ActiveWorkbook

# --- Insert a new "Source" sheet before the existing "Daily" sheet -----
# Final tab order must be: Source, Daily, Hourly
$dailySheet = $wb.Worksheets.Item("Daily")
$sourceSheet = $wb.Worksheets.Add($dailySheet)
$sourceSheet.Name = "Source"

# --- Column width / row sizing for the new sheet -----------------------
$sourceSheet.Columns("B").ColumnWidth = 84.6
$sourceSheet.Rows(2).RowHeight = 20.25

# --- Title text in B2 ----------------------------------------------------
$titleCell = $sourceSheet.Range("B2")
$titleCell.Value = "Descriptions for the relevant columns from the Fitabase data dictionary avaialble at this link:`n"
$titleCell.Font.ThemeColor = 1

# --- Hyperlink text in B3 --------------------------------------------------
$linkUrl = "https://www.fitabase.com/media/2126/fitabase-fitbit-data-dictionary-as-of-05162025.pdf"
$linkCell = $sourceSheet.Range("B3")
$linkCell.Value = $linkUrl
$sourceSheet.Hyperlinks.Add($linkCell, $linkUrl)
$linkCell.Font.Underline = $true
$linkCell.Font.Color = 255 * 65536
